$wb = $excel.ActiveWorkbook

# --- Step 1: write the new raw signal values into Step1_Data (sheet 1) ---
$ws1 = $wb.Worksheets.Item("Step1_Data")

$newValues = @(
    @(2,13,"0.1377531619801652"),
    @(2,14,"0.0002319426830787245"),
    @(2,15,"0.3042291656524011"),
    @(2,16,"0.003411002167705636"),
    @(2,17,"0.01780031285624174"),
    @(2,18,"0.002511147539177135"),
    @(2,19,"0.001096633862132395"),
    @(2,20,"0.002733430917574621"),
    @(2,21,"0.04410713824180821"),
    @(2,22,"0.01538253435154046"),
    @(2,23,"0.0993400432223335"),
    @(2,24,"0.001733258177752026"),
    @(2,25,"0.01875636375177725"),
    @(2,26,"0.008361527981705277"),
    @(2,27,"0.06004233118657825"),
    @(2,28,"0.05444034569682233"),
    @(2,29,"0.01561922586345793"),
    @(2,30,"0.008095957433879863"),
    @(2,31,"0.04481264246228358"),
    @(2,32,"0.005082217371291727"),
    @(2,33,"0.01728164697008704"),
    @(2,34,"0.07005221304364712"),
    @(2,35,"0.03479654689933541"),
    @(2,36,"0.00502815030120207"),
    @(2,37,"0.001886900433956708"),
    @(2,38,"0.000750581607463658"),
    @(2,39,"0.006245293591015495"),
    @(2,40,"0.004526374871873207"),
    @(2,41,"5.287185195168939E-08"),
    @(2,42,"0.005293075167521324"),
    @(2,43,"0.00635090392824317"),
    @(2,44,"0.00224787691409595"),
    @(3,13,"0.08457250662181431"),
    @(3,14,"0.009719235059112647"),
    @(3,15,"0.2871496132964586"),
    @(3,16,"0.006844390927714629"),
    @(3,17,"0.005955146104945607"),
    @(3,18,"0.001155977951904309"),
    @(3,19,"0.0008186382212923966"),
    @(3,20,"0.00123632590379427"),
    @(3,21,"0.04446902074670744"),
    @(3,22,"0.0005110366385942267"),
    @(3,23,"0.09991021107940756"),
    @(3,24,"0.02533365155180051"),
    @(3,25,"0.006170137912919837"),
    @(3,26,"0.008127180105768187"),
    @(3,27,"0.02257921320634892"),
    @(3,28,"0.07600066980982052"),
    @(3,29,"0.0003622031928065767"),
    @(3,30,"0.03163221939319789"),
    @(3,31,"0.0372539431170129"),
    @(3,32,"0.004192430527455261"),
    @(3,33,"0.06126066207827222"),
    @(3,34,"0.02775410632781265"),
    @(3,35,"0.1098900240097356"),
    @(3,36,"0.002606664000882658"),
    @(3,37,"0.01915609277378155"),
    @(3,38,"0.0007082683104214916"),
    @(3,39,"0.0006857765348885152"),
    @(3,40,"0.01242101466806102"),
    @(3,41,"0.0004830368848384983"),
    @(3,42,"0.0007805529357986644"),
    @(3,43,"0.005657522861349884"),
    @(3,44,"0.004602527245280466"),
    @(4,13,"0.185913813103198"),
    @(4,14,"0.05271207212555132"),
    @(4,15,"0.1948833354195419"),
    @(4,16,"0.003062509198281188"),
    @(4,17,"0.009135489992745232"),
    @(4,18,"1.499150404065134E-05"),
    @(4,19,"0.0005644004144444916"),
    @(4,20,"0.001472430918410208"),
    @(4,21,"0.04863592675824966"),
    @(4,22,"0.02325983366860567"),
    @(4,23,"0.08583053681652869"),
    @(4,24,"0.00252387920427389"),
    @(4,25,"0.01306009005091993"),
    @(4,26,"0.002108435638860065"),
    @(4,27,"0.05928897672838954"),
    @(4,28,"0.03642368855150675"),
    @(4,29,"0.02560156718281032"),
    @(4,30,"0.003162121559286415"),
    @(4,31,"0.05555317492915975"),
    @(4,32,"0.01522757745543032"),
    @(4,33,"0.01695819249539269"),
    @(4,34,"0.09404947749922746"),
    @(4,35,"0.02667230162260687"),
    @(4,36,"0.01484441826113343"),
    @(4,37,"0.00367635792633211"),
    @(4,38,"0.001874769955872187"),
    @(4,39,"0.009260870979330534"),
    @(4,40,"0.003931986033282021"),
    @(4,41,"4.352582633482202E-06"),
    @(4,42,"0.002992333874673776"),
    @(4,43,"0.006286292383533729"),
    @(4,44,"0.001013795165747966"),
    @(5,13,"0.1089376336892066"),
    @(5,14,"0.003460973383102741"),
    @(5,15,"0.2920876544476066"),
    @(5,16,"0.005455383103800835"),
    @(5,17,"0.0187911846593974"),
    @(5,18,"0.003167175426139229"),
    @(5,19,"5.201449235441082E-10"),
    @(5,20,"2.417959712673642E-05"),
    @(5,21,"0.03826605877398412"),
    @(5,22,"0.002630477573356998"),
    @(5,23,"0.1153782310525888"),
    @(5,24,"0.006656496448318329"),
    @(5,25,"0.008484203085197824"),
    @(5,26,"0.0148391604365933"),
    @(5,27,"0.03074055563858205"),
    @(5,28,"0.08058781326737044"),
    @(5,29,"0.006204950163751749"),
    @(5,30,"0.02447021839527631"),
    @(5,31,"0.05170999622999889"),
    @(5,32,"2.325886674201987E-06"),
    @(5,33,"0.03897300481327341"),
    @(5,34,"0.05148556741818553"),
    @(5,35,"0.06424474972351045"),
    @(5,36,"0.004850958775573553"),
    @(5,37,"0.006340584338381911"),
    @(5,38,"0.001270266575318844"),
    @(5,39,"0.001982054313488899"),
    @(5,40,"0.007897120453517731"),
    @(5,41,"7.314837059413304E-05"),
    @(5,42,"0.002195725845056544"),
    @(5,43,"0.006191986422027673"),
    @(5,44,"0.00260016117285337"),
    @(6,13,"0.1786757011141007"),
    @(6,14,"0.1654454902998969"),
    @(6,15,"0.1224515099126162"),
    @(6,16,"0.01915912245858102"),
    @(6,17,"0.01863768331423318"),
    @(6,18,"0.0005079020413086957"),
    @(6,19,"0.002046549942462426"),
    @(6,20,"0.00867493729260394"),
    @(6,21,"0.01127093349295471"),
    @(6,22,"0.09417741700176747"),
    @(6,23,"0.02592704050565183"),
    @(6,24,"0.000832519194703629"),
    @(6,25,"0.0322085347712561"),
    @(6,26,"0.002047968375946144"),
    @(6,27,"0.09862888140898696"),
    @(6,28,"0.003698953566135834"),
    @(6,29,"0.04810432976197451"),
    @(6,30,"0.008924893826378453"),
    @(6,31,"0.01804910698303881"),
    @(6,32,"0.03110421237440046"),
    @(6,33,"0.004579049175052191"),
    @(6,34,"0.07874711259148717"),
    @(6,35,"0.0006322414137134375"),
    @(6,36,"0.003355787207144173"),
    @(6,37,"2.272943751283915E-06"),
    @(6,38,"0.0003312660126691367"),
    @(6,39,"0.01184254774644472"),
    @(6,40,"7.169007496951954E-05"),
    @(6,41,"0.0002712854167135838"),
    @(6,42,"0.006329735716960435"),
    @(6,43,"0.00261231949255713"),
    @(6,44,"0.0006510045695393407"),
    @(7,12,"0.02837503411333008"),
    @(7,13,"0.06821546778835849"),
    @(7,14,"0.1420588935183306"),
    @(7,15,"0.0229669163229274"),
    @(7,16,"0.0002668793046895524"),
    @(7,17,"0.0002590190178715127"),
    @(7,18,"1.404916381074246E-05"),
    @(7,19,"0.0001151022180865232"),
    @(7,20,"0.017748896438373"),
    @(7,21,"0.0006345937343309589"),
    @(7,22,"0.1195004689356782"),
    @(7,23,"0.03529554681480033"),
    @(7,24,"0.009865270984007097"),
    @(7,25,"0.04209741026282626"),
    @(7,26,"0.04134531497794237"),
    @(7,27,"0.05923164921459659"),
    @(7,28,"0.000923741566822898"),
    @(7,29,"0.02829887550543361"),
    @(7,30,"0.03349208932276394"),
    @(7,31,"0.04518252461577946"),
    @(7,32,"0.0593334644868666"),
    @(7,33,"0.0004085832397945892"),
    @(7,34,"0.1524366990538717"),
    @(7,35,"0.01774458765103963"),
    @(7,36,"0.01804694477465443"),
    @(7,37,"0.002389070253101914"),
    @(7,38,"0.003762236436606363"),
    @(7,39,"0.008719422062613657"),
    @(7,40,"0.02398157164351821"),
    @(7,41,"7.140692186650437E-05"),
    @(7,42,"0.00183483261245889"),
    @(7,43,"0.01538343704284808"),
    @(8,4,"0.05544700688376016"),
    @(8,5,"0.03713221100471886"),
    @(8,6,"0.1596299274332655"),
    @(8,7,"0.01378162286898858"),
    @(8,8,"6.280007549245903E-05"),
    @(8,12,"0.01714306498549595"),
    @(8,14,"0.114083750970657"),
    @(8,15,"0.02920179206540669"),
    @(8,16,"0.006140765192337181"),
    @(8,17,"0.03428083082759018"),
    @(8,18,"0.05704070273249308"),
    @(8,19,"0.05458042473306978"),
    @(8,20,"0.003105949344837635"),
    @(8,21,"0.02122672879795038"),
    @(8,22,"0.03637641063656266"),
    @(8,23,"0.04004017999128388"),
    @(8,24,"0.07327343988900127"),
    @(8,25,"3.266421104949425E-05"),
    @(8,26,"0.1560680369453143"),
    @(8,27,"0.007174359163032516"),
    @(8,28,"0.02151009528478876"),
    @(8,29,"0.001432279634245869"),
    @(8,30,"0.002552306720021762"),
    @(8,31,"0.01259374744821318"),
    @(8,32,"0.02705159835239276"),
    @(8,34,"0.003022982498011318"),
    @(8,35,"0.0160143213100187"),
    @(9,13,"0.1333231199823689"),
    @(9,14,"0.007805562050949928"),
    @(9,15,"0.1862518881955146"),
    @(9,16,"0.002299009586706642"),
    @(9,17,"0.002015816325200296"),
    @(9,18,"0.004989350807255975"),
    @(9,19,"5.796283024742135E-05"),
    @(9,20,"0.0005186347865228135"),
    @(9,21,"0.01524982622614071"),
    @(9,22,"0.014844696956972"),
    @(9,23,"0.09663509696221653"),
    @(9,24,"0.002675832474871234"),
    @(9,25,"0.004727228026202425"),
    @(9,26,"0.02202748749929933"),
    @(9,27,"0.09926417767680376"),
    @(9,28,"0.02741556292812056"),
    @(9,29,"0.05454493339337255"),
    @(9,30,"0.0007778716726089539"),
    @(9,31,"0.06989888223285785"),
    @(9,32,"0.003122711586979333"),
    @(9,33,"0.04556419581444172"),
    @(9,34,"0.06083913417032001"),
    @(9,35,"0.0794014706672407"),
    @(9,36,"0.01293536347986436"),
    @(9,37,"0.007375627699283004"),
    @(9,38,"0.006642928479082629"),
    @(9,39,"0.001976151378106015"),
    @(9,40,"0.01753917684267948"),
    @(9,41,"0.002841278429046427"),
    @(9,42,"0.0001223203621531309"),
    @(9,43,"0.0122571977984194"),
    @(9,44,"0.00405950267815123"),
    @(10,12,"0.03105802702806259"),
    @(10,13,"0.06928972581844355"),
    @(10,14,"0.2064938094325463"),
    @(10,15,"0.02364763623188139"),
    @(10,16,"0.009921931655374156"),
    @(10,17,"0.0005891700114057567"),
    @(10,18,"0.01360397760679502"),
    @(10,19,"0.0054518666457588"),
    @(10,20,"0.009274929339491363"),
    @(10,21,"0.006203190069465357"),
    @(10,22,"0.1228523065154644"),
    @(10,23,"0.02365816131301969"),
    @(10,24,"9.600105482775519E-06"),
    @(10,25,"0.03659043096418986"),
    @(10,26,"0.02758103459069596"),
    @(10,27,"0.06411552141700695"),
    @(10,28,"0.001829656174672793"),
    @(10,29,"0.03406597628182217"),
    @(10,30,"0.04537834638589372"),
    @(10,31,"0.02615560472520944"),
    @(10,32,"0.0655220620218311"),
    @(10,33,"0.004141932959268971"),
    @(10,34,"0.111594891807817"),
    @(10,35,"0.0005854983415077766"),
    @(10,36,"0.01774738608356543"),
    @(10,37,"0.0001137858998828708"),
    @(10,38,"0.002162907899612905"),
    @(10,39,"0.01466900977905226"),
    @(10,40,"0.007830217582677825"),
    @(10,41,"4.637355048410476E-05"),
    @(10,42,"0.004750038616207948"),
    @(10,43,"0.01306499314541"),
    @(11,12,"0.06631965444707476"),
    @(11,13,"0.0173401467028428"),
    @(11,14,"0.1884058892405248"),
    @(11,15,"0.007684901277072531"),
    @(11,16,"1.46351947491738E-05"),
    @(11,17,"8.560478920743717E-05"),
    @(11,18,"0.002376379876115347"),
    @(11,19,"0.0003155978879655446"),
    @(11,20,"0.0152996003634207"),
    @(11,21,"1.913913644921436E-06"),
    @(11,22,"0.1178107922535183"),
    @(11,23,"0.02059857525447349"),
    @(11,24,"0.001392025444372978"),
    @(11,25,"0.03043526397587486"),
    @(11,26,"0.06276139155937834"),
    @(11,27,"0.05336399937325492"),
    @(11,28,"0.008290901058851044"),
    @(11,29,"0.01414382674507906"),
    @(11,30,"0.04005762768531888"),
    @(11,31,"0.02620092838840313"),
    @(11,32,"0.08539200688947148"),
    @(11,33,"0.003912127011701631"),
    @(11,34,"0.1521835467277378"),
    @(11,35,"0.0001505611651220034"),
    @(11,36,"0.02237197521297403"),
    @(11,37,"9.888831186681814E-05"),
    @(11,38,"0.001265450656362986"),
    @(11,39,"0.01936803200624526"),
    @(11,40,"0.02033701688687929"),
    @(11,41,"3.796636176807554E-05"),
    @(11,42,"0.005259917961480967"),
    @(11,43,"0.01672285537724682")
)

foreach ($item in $newValues) {
    $r = $item[0]
    $c = $item[1]
    $v = [double]$item[2]
    $ws1.Cells.Item($r, $c).Value = $v
}

# --- Step 2: recompute Step2_Sj as the running cumulative sum of Step1_Data, columns B..AR (2..44) ---
$ws2 = $wb.Worksheets.Item("Step2_Sj")

for ($r = 2; $r -le 11; $r++) {
    $cum = 0.0
    for ($c = 2; $c -le 44; $c++) {
        $raw = [double]($ws1.Cells.Item($r, $c).Value())
        $cum = $cum + $raw
        $ws2.Cells.Item($r, $c).Value = $cum
    }
}

# --- Step 3: recompute the Step3_DataPts_* sheets' D/F/G columns from the refreshed Step2_Sj cumulative values ---
$step3Sheets = @("Step3_DataPts_0.5", "Step3_DataPts_0.7", "Step3_DataPts_0.8", "Step3_DataPts_0.9")

foreach ($sheetName in $step3Sheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($r = 2; $r -le 11; $r++) {
        $threshold = [double]($ws.Cells.Item($r, 2).Value())   # column B
        $startIdx  = [double]($ws.Cells.Item($r, 3).Value())   # column C (unchanged)

        $foundCol = -1
        $foundVal = 0.0
        for ($c = 2; $c -le 44; $c++) {
            $cum = [double]($ws2.Cells.Item($r, $c).Value())
            if (($foundCol -eq -1) -and ($cum -ge $threshold)) {
                $foundCol = $c
                $foundVal = $cum
            }
        }

        $pointExceedsIndex = $foundCol - 1
        $ws.Cells.Item($r, 4).Value = $pointExceedsIndex                # column D
        $ws.Cells.Item($r, 6).Value = $foundVal                         # column F
        $ws.Cells.Item($r, 7).Value = $pointExceedsIndex - $startIdx    # column G
    }
}

Write-Output "done"
